$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 61 (weekly update), pushing the
# previously-existing rows 61-83 down to 62-84.
$ws.Rows.Item(61).Insert()

# Populate the newly-inserted row 61 with the new weekly data point.
$ws.Range("A61").Value = 6
$ws.Range("B61").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C61").Value = "Metropolitana"
$ws.Range("D61").Value = 45097
$ws.Range("E61").Value = 13
$ws.Range("F61").Value = 100112035
$ws.Range("G61").Value = "Bruselas (repollito)"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 150
$ws.Range("K61").Value = 16000
$ws.Range("L61").Value = 17000
$ws.Range("M61").Value = 16400
$ws.Range("N61").Value = '$/malla 15 kilos'
$ws.Range("O61").Value = "Provincia de Quillota"
$ws.Range("P61").Value = 1093
$ws.Range("Q61").Value = 15
$ws.Range("R61").Value = "Hortaliza"
